$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.521825
$ws.Range("H2").Value = 9.04365
$ws.Range("I2").Value = 0.05175904312666389
$ws.Range("J2").Value = 0.03566291435357289
$ws.Range("M2").Value = 0.3000656666666667
$ws.Range("N2").Value = 0.900197
$ws.Range("Q2").Value = 1.356844433175
$ws.Range("R2").Value = 8.141066599049999
$ws.Range("S2").Value = 0.05175904312666389
$ws.Range("T2").Value = 0.03566291435357289

$ws.Range("I3").Value = 0.01890340423255524
$ws.Range("J3").Value = 0.019537179749058
$ws.Range("M3").Value = 0.3000656666666667
$ws.Range("N3").Value = 0.900197
$ws.Range("Q3").Value = 0.4955458457420001
$ws.Range("R3").Value = 4.459912611678001
$ws.Range("S3").Value = 0.01890340423255524
$ws.Range("T3").Value = 0.019537179749058

$ws.Range("G4").Value = 14.04838233333333
$ws.Range("H4").Value = 42.145147
$ws.Range("I4").Value = 0.1608047253157438
$ws.Range("J4").Value = 0.166196034552392
$ws.Range("M4").Value = 0.3000656666666667
$ws.Range("N4").Value = 0.900197
$ws.Range("Q4").Value = 4.215437210439889
$ws.Range("R4").Value = 37.938934893959
$ws.Range("S4").Value = 0.1608047253157438
$ws.Range("T4").Value = 0.166196034552392

$ws.Range("G5").Value = 3.9801995
$ws.Range("H5").Value = 7.960399000000001
$ws.Range("I5").Value = 0.04555933004334003
$ws.Range("J5").Value = 0.0313912002075785
$ws.Range("M5").Value = 0.3000656666666667
$ws.Range("N5").Value = 0.900197
$ws.Range("Q5").Value = 1.194321216433833
$ws.Range("R5").Value = 7.165927298603001
$ws.Range("S5").Value = 0.04555933004334003
$ws.Range("T5").Value = 0.0313912002075785

$ws.Range("G6").Value = 51.74080633333333
$ws.Range("H6").Value = 155.222419
$ws.Range("I6").Value = 0.5922508337707373
$ws.Range("J6").Value = 0.6121072613990376
$ws.Range("M6").Value = 0.3000656666666667
$ws.Range("N6").Value = 0.900197
$ws.Range("Q6").Value = 15.52563954628256
$ws.Range("R6").Value = 139.730755916543
$ws.Range("S6").Value = 0.5922508337707373
$ws.Range("T6").Value = 0.6121072613990376

$ws.Range("G7").Value = 11.42032333333333
$ws.Range("H7").Value = 34.26097
$ws.Range("I7").Value = 0.1307226635109598
$ws.Range("J7").Value = 0.1351054097383612
$ws.Range("M7").Value = 0.3000656666666667
$ws.Range("N7").Value = 0.900197
$ws.Range("Q7").Value = 3.426846934565556
$ws.Range("R7").Value = 30.84162241109
$ws.Range("S7").Value = 0.1307226635109598
$ws.Range("T7").Value = 0.1351054097383612
